# Add 2022-Q4 data.
#
# Shape of the edit:
#  - A new "2022-Q4" sheet is inserted right after "总计" (so the existing
#    "2022-Q3" / "2022-Q2" / "2022-Q1" sheets shift one position to the right).
#  - The new sheet's layout is identical to the (old) "2022-Q3" sheet, so we
#    duplicate that sheet and then overwrite the numbers that changed.
#  - The "总计" (totals) summary sheet gets a new row for 2022-Q4 inserted
#    right under the header, and the 2022-Q3 row's counts are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the "2022-Q3" sheet to create "2022-Q4" right before it.
# ---------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($sheetQ3)

$sheetQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$sheetQ4.Name = "2022-Q4"

# D:G hold text-formatted numbers (e.g. "3.47") in the source data, not real
# numbers - force text so Excel doesn't silently convert them back to the
# Number type when we overwrite the values below.
$sheetQ4.Range("D2:G3").NumberFormat = "@"

# Update the Q4 numbers (only D/E/F/G/H change; A/B/C stay the same funds).
$sheetQ4.Range("D2").Value = "3.47"
$sheetQ4.Range("E2").Value = "70.57"
$sheetQ4.Range("F2").Value = "2.49"
$sheetQ4.Range("G2").Value = "0.0864"
$sheetQ4.Range("H2").Value = 9

$sheetQ4.Range("D3").Value = "0.90"
$sheetQ4.Range("E3").Value = "70.57"
$sheetQ4.Range("F3").Value = "2.49"
$sheetQ4.Range("G3").Value = "0.0224"
$sheetQ4.Range("H3").Value = 9

# Drop the temporary "@" number format now that the text values are locked
# in - keeps D2:G3 styled the same plain/default way as the rest of the copied sheet.
$sheetQ4.Range("D2:G3").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a 2022-Q4 row, refresh 2022-Q3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room: copy the last data row down one row (for formatting), then
# fill in the new bottom row (2022-Q1, unchanged values) and shift the
# rest up-to-down by writing new values row by row.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.11

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.29

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.09
